$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I43").Value = 6242.7144
$ws.Range("J43").Value = 5749.857
$ws.Range("K43").Value = 6242.7144
$ws.Range("L43").Value = 5749.857
$ws.Range("M43").Value = -6173.7144
$ws.Range("N43").Value = -5887.857

$ws.Range("H98").Value = 1000.57574
$ws.Range("I98").Value = 563.125
$ws.Range("J98").Value = 14999
$ws.Range("K98").Value = 563.125
$ws.Range("L98").Value = 14999
$ws.Range("M98").Value = 934.875
$ws.Range("N98").Value = -17995

$ws.Range("H111").Value = 1991.25
$ws.Range("J111").Value = 1735
$ws.Range("L111").Value = 5205
$ws.Range("N111").Value = -11339

$ws.Range("H112").Value = 4561.4424
$ws.Range("J112").Value = 4677.24
$ws.Range("L112").Value = 14031.72
$ws.Range("N112").Value = -16247.72

$ws.Range("H122").Value = 1000.57574
$ws.Range("I122").Value = 563.125
$ws.Range("J122").Value = 14999
$ws.Range("K122").Value = 1689.375
$ws.Range("L122").Value = 44997
$ws.Range("M122").Value = 760.625
$ws.Range("N122").Value = -49897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14808.247
$ws.Range("I32").Value = 14488.421
$ws.Range("K32").Value = 14488.421
$ws.Range("M32").Value = -14201.421

$ws.Range("H61").Value = 2854.9534
$ws.Range("I61").Value = 2069.5676
$ws.Range("J61").Value = 7698.1665
$ws.Range("K61").Value = 2069.5676
$ws.Range("L61").Value = 7698.1665
$ws.Range("M61").Value = -1857.5676
$ws.Range("N61").Value = -8122.1665

$ws.Range("H63").Value = 4361.875
$ws.Range("I63").Value = 3865.5557
$ws.Range("K63").Value = 3865.5557
$ws.Range("M63").Value = -3179.5557

$ws.Range("H66").Value = 4361.875
$ws.Range("I66").Value = 3865.5557
$ws.Range("K66").Value = 19327.7785
$ws.Range("M66").Value = -15895.7785

$ws.Range("H128").Value = 84596.8
$ws.Range("J128").Value = 84596.8
$ws.Range("L128").Value = 84596.8
$ws.Range("N128").Value = -94556.8

$ws.Range("H132").Value = 2259.6
$ws.Range("I132").Value = 2259.6
$ws.Range("K132").Value = 6778.799999999999
$ws.Range("M132").Value = -4248.799999999999

$ws.Range("H136").Value = 2854.9534
$ws.Range("I136").Value = 2069.5676
$ws.Range("J136").Value = 7698.1665
$ws.Range("K136").Value = 6208.702799999999
$ws.Range("L136").Value = 23094.4995
$ws.Range("M136").Value = -3658.702799999999
$ws.Range("N136").Value = -28194.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 31426.857
$ws.Range("J75").Value = 34998
$ws.Range("L75").Value = 34998
$ws.Range("N75").Value = -36870

$ws.Range("H78").Value = 31426.857
$ws.Range("J78").Value = 34998
$ws.Range("L78").Value = 104994
$ws.Range("N78").Value = -114354

$ws.Range("H107").Value = 39366.93
$ws.Range("I107").Value = 48558.09
$ws.Range("J107").Value = 5666
$ws.Range("K107").Value = 48558.09
$ws.Range("L107").Value = 5666
$ws.Range("M107").Value = -46638.09
$ws.Range("N107").Value = -9506

$ws.Range("H132").Value = 94223.60000000001
$ws.Range("J132").Value = 94223.60000000001
$ws.Range("L132").Value = 94223.60000000001
$ws.Range("N132").Value = -104343.6

$ws.Range("H134").Value = 2513.0217
$ws.Range("I134").Value = 2348
$ws.Range("J134").Value = 4245.75
$ws.Range("K134").Value = 7044
$ws.Range("L134").Value = 12737.25
$ws.Range("M134").Value = -4509
$ws.Range("N134").Value = -17807.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2972.1052
$ws.Range("I58").Value = 3649.2222
$ws.Range("J58").Value = 1310.091
$ws.Range("K58").Value = 3649.2222
$ws.Range("L58").Value = 1310.091
$ws.Range("M58").Value = -3446.2222
$ws.Range("N58").Value = -1716.091

$ws.Range("H94").Value = 1374.5555
$ws.Range("J94").Value = 1072
$ws.Range("L94").Value = 1072
$ws.Range("N94").Value = -1974

$ws.Range("H132").Value = 28746.863
$ws.Range("I132").Value = 38858.688
$ws.Range("J132").Value = 1782
$ws.Range("K132").Value = 116576.064
$ws.Range("L132").Value = 5346
$ws.Range("M132").Value = -114046.064
$ws.Range("N132").Value = -10406

$ws.Range("H136").Value = 2972.1052
$ws.Range("I136").Value = 3649.2222
$ws.Range("J136").Value = 1310.091
$ws.Range("K136").Value = 10947.6666
$ws.Range("L136").Value = 3930.273
$ws.Range("M136").Value = -8397.6666
$ws.Range("N136").Value = -9030.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1192.7667
$ws.Range("J5").Value = 2289.25
$ws.Range("L5").Value = 6867.75
$ws.Range("N5").Value = -7091.75

$ws.Range("H109").Value = 3165
$ws.Range("I109").Value = 2279.4443
$ws.Range("K109").Value = 6838.3329
$ws.Range("M109").Value = -5798.3329

$ws.Range("H113").Value = 2141.875
$ws.Range("J113").Value = 2260
$ws.Range("L113").Value = 6780
$ws.Range("N113").Value = -11120

$ws.Range("H122").Value = 137971
$ws.Range("I122").Value = 429462.28
$ws.Range("K122").Value = 3865160.52
$ws.Range("M122").Value = -3862710.52

$ws.Range("H129").Value = 2704.2727
$ws.Range("I129").Value = 1324.1818
$ws.Range("J129").Value = 4084.3635
$ws.Range("K129").Value = 3972.5454
$ws.Range("L129").Value = 12253.0905
$ws.Range("M129").Value = 1027.4546
$ws.Range("N129").Value = -22253.0905

$ws.Range("H135").Value = 1192.7667
$ws.Range("J135").Value = 2289.25
$ws.Range("L135").Value = 20603.25
$ws.Range("N135").Value = -25673.25

$ws.Range("H137").Value = 3304.6667
$ws.Range("I137").Value = 2666
$ws.Range("J137").Value = 3943.3333
$ws.Range("K137").Value = 7998
$ws.Range("L137").Value = 11829.9999
$ws.Range("M137").Value = -2898
$ws.Range("N137").Value = -22029.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 57500
$ws.Range("J11").Value = 57500
$ws.Range("L11").Value = 57500
$ws.Range("N11").Value = -57778

$ws.Range("H80").Value = 8027
$ws.Range("I80").Value = 4049.4443
$ws.Range("K80").Value = 4049.4443
$ws.Range("M80").Value = -3051.4443

$ws.Range("H83").Value = 8027
$ws.Range("I83").Value = 4049.4443
$ws.Range("K83").Value = 20247.2215
$ws.Range("M83").Value = -15255.2215

$ws.Range("H113").Value = 2771.1428
$ws.Range("I113").Value = 2771.1428
$ws.Range("K113").Value = 2771.1428
$ws.Range("M113").Value = -601.1428000000001

$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4970
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H61").Value = 1331.3
$ws.Range("J61").Value = 1337.25
$ws.Range("L61").Value = 1337.25
$ws.Range("N61").Value = -1741.25

$ws.Range("H113").Value = 1331.3
$ws.Range("J113").Value = 1337.25
$ws.Range("L113").Value = 1337.25
$ws.Range("N113").Value = -5677.25

$ws.Range("H135").Value = 53250
$ws.Range("J135").Value = 53250
$ws.Range("L135").Value = 53250
$ws.Range("N135").Value = -63390

$ws.Range("H136").Value = 3130.238
$ws.Range("J136").Value = 3416.05
$ws.Range("L136").Value = 10248.15
$ws.Range("N136").Value = -15348.15

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 428.53333
$ws.Range("I100").Value = 430.57144
$ws.Range("K100").Value = 861.14288
$ws.Range("M100").Value = -320.14288

$ws.Range("H113").Value = 823.9
$ws.Range("I113").Value = 589.1429000000001
$ws.Range("K113").Value = 1767.4287
$ws.Range("M113").Value = 402.5712999999998

$ws.Range("H122").Value = 20054924
$ws.Range("I122").Value = 35809228
$ws.Range("K122").Value = 107427684
$ws.Range("M122").Value = -107425234

$ws.Range("H132").Value = 5017252
$ws.Range("I132").Value = 6118149.5
$ws.Range("J132").Value = 2053.111
$ws.Range("K132").Value = 18354448.5
$ws.Range("L132").Value = 6159.333
$ws.Range("M132").Value = -18351918.5
$ws.Range("N132").Value = -11219.333

$ws.Range("H136").Value = 17525.922
$ws.Range("I136").Value = 19960.227
$ws.Range("K136").Value = 59880.681
$ws.Range("M136").Value = -57330.681
